# "script updates regarding clean up and close"
# - devices: swap the test device from Android/Galaxy S7 Edge to iOS/iPhone-6,
#   and point the bundleId/appPackage at the Holder (psm) app instead of the
#   Checker app.
# - signIn: refresh the UAT holder credentials/licence/address and add new
#   credit-card + app build columns used by the updated script.
# - checkerSignIn: add a buildName column for the checker app.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# devices
# ---------------------------------------------------------------------------
$wsDevices = $wb.Worksheets.Item("devices")

$wsDevices.Range("B2").ClearFormats()
$wsDevices.Range("B2").Value = "iOS"

$wsDevices.Range("E2").ClearFormats()
$wsDevices.Range("E2").Value = "iPhone-6"

$wsDevices.Range("L2").Value = "au.gov.nsw.onegov.MyLicences.uat"
$wsDevices.Range("M2").Value = "au.gov.nsw.onegov.app.holder.psm"

$wsDevices.Columns("L:M").AutoFit()

# ---------------------------------------------------------------------------
# signIn
# ---------------------------------------------------------------------------
$wsSignIn = $wb.Worksheets.Item("signIn")

$wsSignIn.Range("A2").Value = "automation_psm02@yopmail.com"

$wsSignIn.Range("D2").Value = "'RS2694206"

$wsSignIn.Range("E2").ClearFormats()
$wsSignIn.Range("E2").Value = "'08 Aug 2013"

$wsSignIn.Range("F2").ClearFormats()
$wsSignIn.Range("F2").Value = "'08 Aug 2017"

$wsSignIn.Range("K2").Value = "2-24 Rawson Pl,HAYMARKET,NSW 2000"
$wsSignIn.Range("L2").Value = "Jacky Williams"

# new header columns (M:Q credit-card fields, R appBuildName, S appVersion)
$wsSignIn.Range("M1:S1").Interior.Color = 65535
$wsSignIn.Range("M1").Value = "cardNumber"
$wsSignIn.Range("N1").Value = "cardExpiryMonth"
$wsSignIn.Range("O1").Value = "cardExpiryYear"
$wsSignIn.Range("P1").Value = "cardCVVNum"
$wsSignIn.Range("Q1").Value = "cardName"
$wsSignIn.Range("R1").Value = "'appBuildName"
$wsSignIn.Range("S1").Value = "appVersion"

# new data row for the credit-card + app build columns
$wsSignIn.Range("M2").Value = "'5163200000000008"
$wsSignIn.Range("N2").Value = "'08-Aug"
$wsSignIn.Range("N2").NumberFormat = "d-mmm"
$wsSignIn.Range("O2").Value = "'2020"
$wsSignIn.Range("P2").Value = "'070"
$wsSignIn.Range("Q2").Value = "Srikanth"
$wsSignIn.Range("R2").Value = "UAT Holder"
$wsSignIn.Range("S2").Value = "'1.0.2-PSM"

$wsSignIn.Range("F9").Select() | Out-Null

# ---------------------------------------------------------------------------
# checkerSignIn
# ---------------------------------------------------------------------------
$wsChecker = $wb.Worksheets.Item("checkerSignIn")

$wsChecker.Range("M1").Interior.Color = 65535
$wsChecker.Range("M1").Value = "buildName"
$wsChecker.Range("M2").Value = "UAT Checker"

$wsChecker.Range("G9").Select() | Out-Null

$wsDevices.Select() | Out-Null
